# The document has a "_GoBack" bookmark (the position Word restores the
# cursor to on reopen) sitting after the spaces that follow "Секретарь".
# This edit moves that bookmark so it instead sits at the very start of
# the document (immediately before its first run of text), and removes
# it from its old location.
#
# Word's "_GoBack" bookmark is a singleton: (re)adding a bookmark named
# "_GoBack" anywhere automatically replaces/removes any existing one with
# that name, so a single Bookmarks.Add(...) call at the new location takes
# care of both the "add at top" and "remove from old spot" halves of the
# diff in one step.
#
# Inserting a *collapsed* (zero-length) bookmark exactly at document
# position 0 needs a small workaround: insert a one-character placeholder
# at the very start, wrap the bookmark around that single character (which
# reliably collapses once we delete it), then delete the placeholder. The
# net text content of the document is unchanged, but we're left with an
# empty bookmarkStart/bookmarkEnd pair sitting right before the first run.

$d = $word.ActiveDocument

$start = $d.Range(0, 0)
$start.InsertBefore("X")

$around = $d.Range(0, 1)
$around.Bookmarks.Add("_GoBack")

$placeholder = $d.Range(0, 1)
$placeholder.Text = ""

Write-Output "Moved _GoBack bookmark to the start of the document."
